$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blockA = New-Object "object[,]" 3,6
$blockA[0,0] = 27.22363662719727
$blockA[0,1] = 31.01375389099121
$blockA[0,2] = 31.3092212677002
$blockA[0,3] = 26.89251136779785
$blockA[0,4] = 60594610
$blockA[0,5] = "VRNT"
$blockA[1,0] = 31.33978652954102
$blockA[1,1] = 32.94447326660156
$blockA[1,2] = 33.65766525268555
$blockA[1,3] = 31.1054515838623
$blockA[1,4] = 60594610
$blockA[1,5] = "VRNT"
$blockA[2,0] = 29.7554759979248
$blockA[2,1] = 27.16760063171387
$blockA[2,2] = 30.40753936767578
$blockA[2,3] = 26.23026084899902
$blockA[2,4] = 60594610
$blockA[2,5] = "VRNT"
$ws.Range("D2:I4").Value = $blockA

$ws.Range("H5").Value = 60594610

$blockB = New-Object "object[,]" 39,6
$blockB[0,0] = 18.5379524230957
$blockB[0,1] = 18.0998477935791
$blockB[0,2] = 19.10341262817383
$blockB[0,3] = 15.16046905517578
$blockB[0,4] = 60594610
$blockB[0,5] = "VRNT"
$blockB[1,0] = 17.3204288482666
$blockB[1,1] = 16.805908203125
$blockB[1,2] = 17.58532905578613
$blockB[1,3] = 16.13856315612793
$blockB[1,4] = 60594610
$blockB[1,5] = "VRNT"
$blockB[2,0] = 17.95720863342285
$blockB[2,1] = 17.38155937194824
$blockB[2,2] = 18.8792667388916
$blockB[2,3] = 17.11156463623047
$blockB[2,4] = 60594610
$blockB[2,5] = "VRNT"
$blockB[3,0] = 18.31380462646484
$blockB[3,1] = 19.12888336181641
$blockB[3,2] = 19.84207725524902
$blockB[3,3] = 17.98267936706543
$blockB[3,4] = 60594610
$blockB[3,5] = "VRNT"
$blockB[4,0] = 19.10341262817383
$blockB[4,1] = 19.23076820373535
$blockB[4,2] = 20.12226104736328
$blockB[4,3] = 18.28833389282227
$blockB[4,4] = 60594610
$blockB[4,5] = "VRNT"
$blockB[5,0] = 20.0967903137207
$blockB[5,1] = 20.93734169006348
$blockB[5,2] = 21.19205284118652
$blockB[5,3] = 19.74019432067871
$blockB[5,4] = 60594610
$blockB[5,5] = "VRNT"
$blockB[6,0] = 20.3005599975586
$blockB[6,1] = 20.22414779663086
$blockB[6,2] = 20.3005599975586
$blockB[6,3] = 18.87417221069336
$blockB[6,4] = 60594610
$blockB[6,5] = "VRNT"
$blockB[7,0] = 21.67600631713867
$blockB[7,1] = 22.28731536865234
$blockB[7,2] = 22.69485473632812
$blockB[7,3] = 20.70809936523437
$blockB[7,4] = 60594610
$blockB[7,5] = "VRNT"
$blockB[8,0] = 21.19205284118652
$blockB[8,1] = 19.81660652160645
$blockB[8,2] = 21.26846694946289
$blockB[8,3] = 18.95058631896973
$blockB[8,4] = 60594610
$blockB[8,5] = "VRNT"
$blockB[9,0] = 21.37035179138184
$blockB[9,1] = 21.49770736694336
$blockB[9,2] = 22.46561431884766
$blockB[9,3] = 20.60621452331543
$blockB[9,4] = 60594610
$blockB[9,5] = "VRNT"
$blockB[10,0] = 22.9495677947998
$blockB[10,1] = 24.73255157470703
$blockB[10,2] = 25.44574546813965
$blockB[10,3] = 22.7203254699707
$blockB[10,4] = 60594610
$blockB[10,5] = "VRNT"
$blockB[11,0] = 23.38767242431641
$blockB[11,1] = 23.14314842224121
$blockB[11,2] = 24.42180252075196
$blockB[11,3] = 21.3652572631836
$blockB[11,4] = 60594610
$blockB[11,5] = "VRNT"
$blockB[12,0] = 24.72236442565918
$blockB[12,1] = 27.12684631347656
$blockB[12,2] = 27.6209888458252
$blockB[12,3] = 24.63576126098633
$blockB[12,4] = 60594610
$blockB[12,5] = "VRNT"
$blockB[13,0] = 30.85073852539062
$blockB[13,1] = 28.90983200073243
$blockB[13,2] = 32.01222610473633
$blockB[13,3] = 27.79419326782227
$blockB[13,4] = 60594610
$blockB[13,5] = "VRNT"
$blockB[14,0] = 29.43453979492188
$blockB[14,1] = 27.14722442626953
$blockB[14,2] = 29.78094863891602
$blockB[14,3] = 26.37289810180664
$blockB[14,4] = 60594610
$blockB[14,5] = "VRNT"
$blockB[15,0] = 23.24503326416016
$blockB[15,1] = 24.76311874389648
$blockB[15,2] = 24.96688652038575
$blockB[15,3] = 23.12786674499512
$blockB[15,4] = 60594610
$blockB[15,5] = "VRNT"
$blockB[16,0] = 29.68415641784668
$blockB[16,1] = 27.95720863342285
$blockB[16,2] = 30.4788589477539
$blockB[16,3] = 26.14365768432617
$blockB[16,4] = 60594610
$blockB[16,5] = "VRNT"
$blockB[17,0] = 21.20224189758301
$blockB[17,1] = 23.62200736999512
$blockB[17,2] = 24.34539031982422
$blockB[17,3] = 20.40753936767578
$blockB[17,4] = 60594610
$blockB[17,5] = "VRNT"
$blockB[18,0] = 22.87315368652344
$blockB[18,1] = 24.22822189331055
$blockB[18,2] = 24.4116153717041
$blockB[18,3] = 21.9001522064209
$blockB[18,4] = 60594610
$blockB[18,5] = "VRNT"
$blockB[19,0] = 24.78858947753906
$blockB[19,1] = 29.01681137084961
$blockB[19,2] = 31.13092231750488
$blockB[19,3] = 24.5236873626709
$blockB[19,4] = 60594610
$blockB[19,5] = "VRNT"
$blockB[20,0] = 38.00815200805664
$blockB[20,1] = 49.29000091552734
$blockB[20,2] = 52.70000076293945
$blockB[20,3] = 37.78400421142578
$blockB[20,4] = 60594610
$blockB[20,5] = "VRNT"
$blockB[21,0] = 48.9900016784668
$blockB[21,1] = 46.11000061035156
$blockB[21,2] = 49.2599983215332
$blockB[21,3] = 43.77999877929688
$blockB[21,4] = 60594610
$blockB[21,5] = "VRNT"
$blockB[22,0] = 42.68999862670898
$blockB[22,1] = 44.63999938964844
$blockB[22,2] = 45.13999938964844
$blockB[22,3] = 41.7400016784668
$blockB[22,4] = 60594610
$blockB[22,5] = "VRNT"
$blockB[23,0] = 46.66999816894531
$blockB[23,1] = 47.59000015258789
$blockB[23,2] = 48.93000030517578
$blockB[23,3] = 46.4900016784668
$blockB[23,4] = 60594610
$blockB[23,5] = "VRNT"
$blockB[24,0] = 51.43000030517578
$blockB[24,1] = 50.22999954223633
$blockB[24,2] = 54.29999923706055
$blockB[24,3] = 48.47999954223633
$blockB[24,4] = 60594610
$blockB[24,5] = "VRNT"
$blockB[25,0] = 54.33000183105469
$blockB[25,1] = 51.04000091552734
$blockB[25,2] = 55.45000076293945
$blockB[25,3] = 47.66999816894531
$blockB[25,4] = 60594610
$blockB[25,5] = "VRNT"
$blockB[26,0] = 45.18999862670898
$blockB[26,1] = 48.4900016784668
$blockB[26,2] = 51.0099983215332
$blockB[26,3] = 44.16999816894531
$blockB[26,4] = 60594610
$blockB[26,5] = "VRNT"
$blockB[27,0] = 35.86999893188477
$blockB[27,1] = 39.33000183105469
$blockB[27,2] = 40.15999984741211
$blockB[27,3] = 32.81000137329102
$blockB[27,4] = 60594610
$blockB[27,5] = "VRNT"
$blockB[28,0] = 37.95000076293945
$blockB[28,1] = 37.38000106811523
$blockB[28,2] = 40.70999908447266
$blockB[28,3] = 36.65999984741211
$blockB[28,4] = 60594610
$blockB[28,5] = "VRNT"
$blockB[29,0] = 36.5099983215332
$blockB[29,1] = 35.88000106811523
$blockB[29,2] = 36.70000076293945
$blockB[29,3] = 33.11000061035156
$blockB[29,4] = 60594610
$blockB[29,5] = "VRNT"
$blockB[30,0] = 37.09999847412109
$blockB[30,1] = 32.38999938964844
$blockB[30,2] = 37.40999984741211
$blockB[30,3] = 31.63999938964844
$blockB[30,4] = 60594610
$blockB[30,5] = "VRNT"
$blockB[31,0] = 18.72999954223633
$blockB[31,1] = 24.56999969482422
$blockB[31,2] = 25.60000038146973
$blockB[31,3] = 18.45000076293945
$blockB[31,4] = 60594610
$blockB[31,5] = "VRNT"
$blockB[32,0] = 29.79000091552734
$blockB[32,1] = 31.61000061035156
$blockB[32,2] = 32.7400016784668
$blockB[32,3] = 28.42000007629395
$blockB[32,4] = 60594610
$blockB[32,5] = "VRNT"
$blockB[33,0] = 30.27000045776367
$blockB[33,1] = 29.65999984741211
$blockB[33,2] = 33.2400016784668
$blockB[33,3] = 28.8799991607666
$blockB[33,4] = 60594610
$blockB[33,5] = "VRNT"
$blockB[34,0] = 36.09000015258789
$blockB[34,1] = 31.54999923706055
$blockB[34,2] = 36.54999923706055
$blockB[34,3] = 30.34000015258789
$blockB[34,4] = 60594610
$blockB[34,5] = "VRNT"
$blockB[35,0] = 21.45999908447266
$blockB[35,1] = 25.20000076293945
$blockB[35,2] = 26.1299991607666
$blockB[35,3] = 21.28000068664551
$blockB[35,4] = 60594610
$blockB[35,5] = "VRNT"
$blockB[36,0] = 24.69000053405762
$blockB[36,1] = 22.56999969482422
$blockB[36,2] = 26.88999938964844
$blockB[36,3] = 22.17000007629395
$blockB[36,4] = 60594610
$blockB[36,5] = "VRNT"
$blockB[37,0] = 18.1299991607666
$blockB[37,1] = 17.54000091552734
$blockB[37,2] = 19.20999908447266
$blockB[37,3] = 16.35000038146973
$blockB[37,4] = 60594610
$blockB[37,5] = "VRNT"
$blockB[38,0] = 20.85000038146973
$blockB[38,1] = 20.38999938964844
$blockB[38,2] = 21.8799991607666
$blockB[38,3] = 18.39999961853028
$blockB[38,4] = 60594610
$blockB[38,5] = "VRNT"
$ws.Range("D6:I44").Value = $blockB
